# Append new trading-log rows (164-167) to Sheet1, mirroring the existing
# TRADING_ATTEMPT / POSITION_FAILED row layout (columns A-L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 164; A = "2026-01-06T12:47:00.521657"; B = "TRADING_ATTEMPT"; C = "ENA";  D = "UNKNOWN"; E = 0.2537161193615907; K = "ATTEMPT"; L = "Attempting trade 1/2" },
    @{ Row = 165; A = "2026-01-06T12:47:00.749178"; B = "POSITION_FAILED"; C = "ENA";  D = "UNKNOWN"; E = $null;               K = "FAILED";  L = "Trade execution failed for trade 1" },
    @{ Row = 166; A = "2026-01-06T12:47:00.801470"; B = "TRADING_ATTEMPT"; C = "DOGE"; D = "UNKNOWN"; E = 0.1519393256021509; K = "ATTEMPT"; L = "Attempting trade 2/2" },
    @{ Row = 167; A = "2026-01-06T12:47:01.125724"; B = "POSITION_FAILED"; C = "DOGE"; D = "UNKNOWN"; E = $null;               K = "FAILED";  L = "Trade execution failed for trade 2" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    if ($null -eq $r.E) {
        $ws.Cells.Item($row, 5).Value = ""
    } else {
        $ws.Cells.Item($row, 5).Value = $r.E
    }
    $ws.Cells.Item($row, 6).Value = ""
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = ""
    $ws.Cells.Item($row, 9).Value = ""
    $ws.Cells.Item($row, 10).Value = ""
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
}
